$d = $word.ActiveDocument

$pairs = @(
    @("2024-04-03 Wednesday", "2024-04-04 Thursday"),
    @("624×8=4992", "501×3=1503"),
    @("161×9=1449", "941×3=2823"),
    @("628×7=4396", "811×9=7299"),
    @("341×8=2728", "657×3=1971"),
    @("980×9=8820", "551×2=1102"),
    @("945×5=4725", "702×5=3510"),
    @("275×8=2200", "204×8=1632"),
    @("881×5=4405", "300×6=1800"),
    @("612×4=2448", "575×6=3450"),
    @("403×9=3627", "862×5=4310"),
    @("914×9=8226", "888×8=7104"),
    @("430×7=3010", "228×8=1824"),
    @("293×4=1172", "195×6=1170"),
    @("743×5=3715", "953×4=3812"),
    @("431×6=2586", "869×8=6952"),
    @("333×9=2997", "787×6=4722"),
    @("918×9=8262", "203×3=609"),
    @("925×8=7400", "264×6=1584"),
    @("690×2=1380", "279×7=1953"),
    @("791×8=6328", "752×8=6016"),
    @("784×3=2352", "898×8=7184"),
    @("548×5=2740", "818×5=4090"),
    @("699×4=2796", "348×9=3132"),
    @("193×3=579", "118×2=236"),
    @("630×2=1260", "846×7=5922")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
